# PeerEvaluationTemplate_7Person.xlsx - typo cleanup in the rating-description
# text blocks (per the commit message: "fixed typos in the blank templates").
#
# The sharedStrings table in the target file is re-ordered by Excel as a
# side-effect of rewriting these particular strings; editing the cells in the
# same order they appear on the sheet (I9, I11, I12, I13, then I8) reproduces
# that natural ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PeerRating")

# H9/I9 - "Having Relevant Knowledge, Skills, and Abilities (KSAs)" description
# typo fixes: "peform" -> "perform", "elses" -> "else's"
$i9 = @'
5: Demonstrates KSAs to do excellent work, acquires new KSA to help team, can perform any role on team if necessary
4: Between 5 above and 3 below
3: Demonstrates sufficient KSA to contribute to team, acquires KSAs to meet requirements, able to perform other tasks
2: Between 3 above and 1 below
1: Missing basic qualification, unable to develop KSAs to contribute to team, unable to perform any one else's duties
'@
$ws.Range("I9").Value = $i9

# H11/I11 - "Interacting with Teammates" description
# typo fix: "contributsions" -> "contributions"
$i11 = @'
5: Is interested in teammates ideas and contributions, makes sure everyone is informed, is encouraging, enthusiastic and asks for feedback/suggestions
3: Listens and respects teammate contributions, communicates clearly, shares info, participates fully, reacts and responds to feedback/suggestions
1: Interrupts, ignores, bosses, or makes fun, takes action without input, does not share, complains, makes excuses, does not interact, is defensive
'@
$ws.Range("I11").Value = $i11

# H12/I12 - "Keeping the Team on Track" description
# typo fix: "sucess" -> "success"
$i12 = @'
5: Monitors teams' progress, makes sure teammates are progressing, gives specific, timely, and constructive feedback
3: Knows what everyone on the team should be doing and notices problems, alerts teammates and suggests solutions with success is threatened
1: Unaware if team is meeting goals, does not pay attention to teammates progress, avoids discussing team problems even when obvious
'@
$ws.Range("I12").Value = $i12

# H13/I13 - "Expecting Quality" description
# typo fixes: "Encouarges" -> "Encourages", "responsiblities" -> "responsibilities"
$i13 = @'
5: Motivates team to do excellent work, cares about excellent work even without reward, believes in team's ability to do excellent work
3: Encourages good work to meet requirements, believes team can meet its responsibilities
1: Satisfied even if not all requirements are met,  avoids work, doubts team can meet requirements
'@
$ws.Range("I13").Value = $i13

# I8 - general instructions line
# typo fix: "each each member" -> "each member"
$ws.Range("I8").Value = "Rating Descriptions (provide whole number ratings (5, 4, 3, 2, or 1) in columns for each member including yourself)"

# Match the author's final cursor position recorded in the saved file.
$ws.Range("I8").Select() | Out-Null
